$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.692.89"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.539.95"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'603.34"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "'144.34"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").Value = "3.539.78"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "'0.404"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").Value = "4.127.48"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -5.60%  "
$ws.Range("D15").Value = "'28.53"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("D16").Value = "3.507.48"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "65.648.49"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'11.03"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'14.30"
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("D22").Value = "'417.48"
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").Value = "'0.599"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").Value = "'77.87"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").Value = "3.677.36"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").Value = "'2.46"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  -3.44%  "
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "3.549.29"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "'0.154"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'24.34"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'1.32"
$ws.Range("E36").Value = "  -9.08%  "
$ws.Range("D37").Value = "'7.52"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("D38").Value = "'174.84"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'5.28"
$ws.Range("E39").Value = "  -6.29%  "
$ws.Range("D40").Value = "'1.58"
$ws.Range("E40").Value = "  -8.18%  "
$ws.Range("D41").Value = "'0.0815"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "'5.07"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "'0.858"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").Value = "'45.04"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "'1.78"
$ws.Range("E45").Value = "  -7.74%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  -6.72%  "
$ws.Range("D48").Value = "'23.52"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").Value = "'1.11"
$ws.Range("E50").Value = "  -7.44%  "
$ws.Range("D51").Value = "'0.906"
$ws.Range("E51").Value = "  -4.18%  "
